$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @("ECs","Il17c","Il17ra","ECs",3,1,0.8163016666666666,2.448905,0.2413881986261334,0.2413881986261334,3,1,1.331285,3.993855,0.0499894697877283,0.04998946978772829,1.086730164308333,9.780571478775,0.01206686806233525,0.01206686806233525),
  @("ECs","Il17c","Il17ra","FAPs",3,1,0.8163016666666666,2.448905,0.2413881986261334,0.2413881986261334,3,1,10.37764133333333,31.132924,0.3896782341125657,0.3896782341125657,8.471285916468888,76.24157324821999,0.09406372697624493,0.09406372697624495),
  @("ECs","Il17c","Il17ra","M2",3,1,0.8163016666666666,2.448905,0.2413881986261334,0.2413881986261334,3,1,12.55252066666667,37.657562,0.4713444924461467,0.4713444924461467,10.24664354106778,92.21979186960999,0.1137769979639245,0.1137769979639245),
  @("ECs","Il17c","Il17ra","sCs",3,1,0.8163016666666666,2.448905,0.2413881986261334,0.2413881986261334,3,1,2.369861666666667,7.109585,0.08898780365355935,0.08898780365355935,1.934522028269444,17.410698254425,0.02148060562362874,0.02148060562362875),
  @("FAPs","Il17c","Il17ra","ECs",3,1,1.210478333333333,3.631435,0.3579500034006598,0.3579500034006599,3,1,1.331285,3.993855,0.0499894697877283,0.04998946978772829,1.611491647991667,14.503424831925,0.01789373088051453,0.01789373088051453),
  @("FAPs","Il17c","Il17ra","FAPs",3,1,1.210478333333333,3.631435,0.3579500034006598,0.3579500034006599,3,1,10.37764133333333,31.132924,0.3896782341125657,0.3896782341125657,12.56190998510444,113.05718986594,0.139485325225756,0.139485325225756),
  @("FAPs","Il17c","Il17ra","M2",3,1,1.210478333333333,3.631435,0.3579500034006598,0.3579500034006599,3,1,12.55252066666667,37.657562,0.4713444924461467,0.4713444924461467,15.19455429571889,136.75098866147,0.1687177626739805,0.1687177626739805),
  @("FAPs","Il17c","Il17ra","sCs",3,1,1.210478333333333,3.631435,0.3579500034006598,0.3579500034006599,3,1,2.369861666666667,7.109585,0.08898780365355935,0.08898780365355935,2.868666200497222,25.817995804475,0.03185318462040881,0.03185318462040882),
  @("M2","Il17c","Il17ra","ECs",3,1,1.210812666666667,3.632438,0.3580488689602556,0.3580488689602557,3,1,1.331285,3.993855,0.0499894697877283,0.04998946978772829,1.611936740943334,14.50743066849,0.01789867311741899,0.01789867311741899),
  @("M2","Il17c","Il17ra","FAPs",3,1,1.210812666666667,3.632438,0.3580488689602556,0.3580488689602557,3,1,10.37764133333333,31.132924,0.3896782341125657,0.3896782341125657,12.56537957652356,113.088416188712,0.1395238509824339,0.1395238509824339),
  @("M2","Il17c","Il17ra","M2",3,1,1.210812666666667,3.632438,0.3580488689602556,0.3580488689602557,3,1,12.55252066666667,37.657562,0.4713444924461467,0.4713444924461467,15.19875102179511,136.788759196156,0.1687643624109886,0.1687643624109886),
  @("M2","Il17c","Il17ra","sCs",3,1,1.210812666666667,3.632438,0.3580488689602556,0.3580488689602557,3,1,2.369861666666667,7.109585,0.08898780365355935,0.08898780365355935,2.869458524247778,25.82512671823,0.03186198244941423,0.03186198244941423),
  @("sCs","Il17c","Il17ra","ECs",2,0.6666666666666666,0.144104,0.432312,0.04261292901295109,0.0426129290129511,3,1,1.331285,3.993855,0.0499894697877283,0.04998946978772829,0.19184349364,1.72659144276,0.00213019772745953,0.00213019772745953),
  @("sCs","Il17c","Il17ra","FAPs",2,0.6666666666666666,0.144104,0.432312,0.04261292901295109,0.0426129290129511,3,1,10.37764133333333,31.132924,0.3896782341125657,0.3896782341125657,1.495459626698667,13.459136640288,0.0166053309281309,0.0166053309281309),
  @("sCs","Il17c","Il17ra","M2",2,0.6666666666666666,0.144104,0.432312,0.04261292901295109,0.0426129290129511,3,1,12.55252066666667,37.657562,0.4713444924461467,0.4713444924461467,1.808868438149333,16.279815943344,0.02008536939725311,0.02008536939725311),
  @("sCs","Il17c","Il17ra","sCs",2,0.6666666666666666,0.144104,0.432312,0.04261292901295109,0.0426129290129511,3,1,2.369861666666667,7.109585,0.08898780365355935,0.08898780365355935,0.3415065456133333,3.07355891052,0.003792030960107554,0.003792030960107555),
)

$startRow = 2
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $r = $startRow + $i
  $rowData = $newRows[$i]
  for ($c = 1; $c -le $rowData.Length; $c++) {
    $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
  }
}
